$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 19: Task 19 Projects Delete (CRUD) now Complete, with a real date value ---
$ws.Range("B14").Copy() | Out-Null
$ws.Range("B19").PasteSpecial(-4122) | Out-Null
$ws.Range("B19").Value = "Complete"
$ws.Range("D19").Value = 45326
$ws.Range("D4").Copy() | Out-Null
$ws.Range("D19").PasteSpecial(-4122) | Out-Null

# --- Row 29: Task 29 gets a full title + Pending status ---
$ws.Range("A29").Value = "Task 29: Frontend Polishing (Projects CRUD)"
$ws.Range("B22").Copy() | Out-Null
$ws.Range("B29").PasteSpecial(-4122) | Out-Null
$ws.Range("B29").Value = "Pending"

# --- Row 30: Task 30 gets a full title + Pending status ---
$ws.Range("A30").Value = "Task 30: Frontend Polishing (Tasks CRUD)"
$ws.Range("B22").Copy() | Out-Null
$ws.Range("B30").PasteSpecial(-4122) | Out-Null
$ws.Range("B30").Value = "Pending"

# --- Row 31: Task 31 gets a full title + Pending status ---
$ws.Range("A31").Value = "Task 31: Frontend Polishing (Members CRUD)"
$ws.Range("B22").Copy() | Out-Null
$ws.Range("B31").PasteSpecial(-4122) | Out-Null
$ws.Range("B31").Value = "Pending"

# --- Row 32: Task 32 label gets a trailing space ---
$ws.Range("A32").Value = "Task 32: "

# --- Style: date format xf (used by D28) gains left alignment ---
$ws.Range("D28").HorizontalAlignment = -4131

# --- Sheet view: clear frozen top-left cell, move selection to E16 ---
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("E16").Select() | Out-Null
